# Apply the "I0 and IF added" edit: add two new columns (I and J) with
# headers "I0" and "IF" in row 1, and numeric values for rows 2-69.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting used by the existing header cells (bold,
# bordered, centered) by copying the format from H1 onto I1:J1.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data values for rows 2-69 ---
$iVals = @(8,8,8,7,8,9,9,8,8,7,7,8,9,9,8,7,8,8,8,8,9,10,11,8,9,9,6,7,8,9,8,9,8,7,8,9,9,9,8,9,8,8,6,8,9,8,8,8,8,8,10,7,7,8,7,7,8,8,5,7,8,4,8,4,9,3,9,9)
$jVals = @(8,8,8,8,8,9,9,8,9,8,8,8,9,9,8,7,9,8,8,9,9,10,11,8,9,9,7,7,8,9,8,9,8,8,8,9,9,9,8,9,8,9,6,8,9,8,8,8,8,8,10,7,7,8,7,7,8,8,6,8,8,5,8,4,9,3,9,9)

for ($i = 0; $i -lt $iVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$i]
    $ws.Cells.Item($row, 10).Value = $jVals[$i]
}

Write-Host "Added columns I (I0) and J (IF) for rows 1-69"
